$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- dbRDA row update: newest results ---
# Update the explanatory footnote first, so shared-string ordering matches
# a natural "edit note, then fill in row" authoring sequence.
$ws.Range("A20").Value = "Note: dbRDA is an ordination technique and does not have p-values. It produced ordination plots in accordance with the nature of the tests (i.e. positive/negative) and appears to produce meaningful results."

# Replace the old "TODO"/"See note" placeholders in the dbRDA row with the
# new verdict now that results are in.
$ws.Range("B10").Value = "Meaningful"
$ws.Range("C10").Value = "Meaningful"
$ws.Range("D10").Value = "Meaningful"
$ws.Range("E10").Value = "Meaningful"
$ws.Range("F10").Value = "Meaningful"
$ws.Range("I10").Value = "Meaningful"
$ws.Range("J10").Value = "Meaningful"

# --- Print / layout tweaks that came with the update ---
$ws.PageSetup.Zoom = 45
$null = $ws.HPageBreaks.Add($ws.Range("A21"))
$null = $ws.VPageBreaks.Add($ws.Range("Q1"))

# --- Misc workbook/view state ---
$wb.CheckCompatibility = $true
$null = $ws.Range("C17").Select()
